$d = $word.ActiveDocument

# Items that move from "to do" (red text) to "done" (green + strikethrough),
# matching the styling already used elsewhere in the document for
# completed items.
$doneTexts = @(
    "Work on Residential Tab",
    "Add in Jessica’s Maps",
    "Type up some intro for the sidebar",
    "Put together some statistical analysis",
    "Lay it out nice for the fourth tab"
)

# wdColor value for RGB 00B050 (Word stores colors as 0xBBGGRR).
$greenColor = 5287936

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.Trim()
    if ($doneTexts -contains $text) {
        $para.Range.Font.StrikeThrough = $true
        $para.Range.Font.Color = $greenColor
    }
}
